# Lecture partielle de l'EDT M1 MIAGE.
# Shift the two recorded schedule dates forward by 3 years (1096 days) and
# update the corresponding day-of-week labels to match the new dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: date 2023-03-17 (Friday, "vendredi") -> 2026-03-17 (Tuesday, "mardi")
$ws.Range("A2").Value = 46098
$ws.Range("B2").Value = "mardi"

# Row 5: date 2023-04-03 (Monday, "lundi") -> 2026-04-03 (Friday, "vendredi")
$ws.Range("A5").Value = 46115
$ws.Range("B5").Value = "vendredi"
